$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.091386431091322606,
    0.042677264653139702,
    -0.0039999999839732681,
    -0.007999999969877436,
    -0.0029999999851328951,
    -0.0019999999862694295,
    -0.0099999999602617962,
    -0.0099999999601134704,
    -0.001999999986445733,
    -0.0019999999868698382,
    -0.0029999999837064806,
    -0.0034999999821212491,
    0.031985897519440698,
    -0.0079999999683488809,
    -0.00099999999087785341,
    -0.0019999999873689944,
    -0.0019999999870519147,
    -0.0039999999804312125,
    -0.0039999999866036084,
    -0.0039999999857425195,
    -0.0039999999855870882,
    -0.0039999999854725132,
    -0.0049999999797165628,
    -0.019999999929118495,
    -0.019999999928212553,
    -0.0024999999843569043,
    -0.0024999999838684062,
    0.05766321521982487,
    -0.0069999999636687349,
    -0.059999999789216218,
    -0.0069999999603265195,
    -0.0099999999503790349,
    -0.0039999999698512312
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
